# Add the two remaining review-comment rows (51 & 52) to the "Comments"
# sheet, plus the small structural housekeeping that Excel performs when
# new rows are inserted into the reviewed range (row 62 / row 130 filler
# rows, dimension, data validation ranges, hyperlinks, selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comments")
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Create rows 51 and 52 by copying the formatting of the previous two
#    data rows (49 & 50) so that styles/row-height match exactly, then
#    overwrite the copied values with the new comment data.
# ---------------------------------------------------------------------
$ws.Range("A49:L49").Copy($ws.Range("A51:L51"))
$ws.Range("A50:L50").Copy($ws.Range("A52:L52"))
$ws.Rows.Item(51).RowHeight = 56.25
$ws.Rows.Item(52).RowHeight = 56.25

# ---- Row 51 --------------------------------------------------------
$ws.Cells.Item(51,1).Value  = "Samuel J. Crawford"
$ws.Cells.Item(51,2).Value  = "crawfs1@mcmaster.ca"
$ws.Cells.Item(51,3).Value  = "McMaster University"
$ws.Cells.Item(51,4).Value  = "Yes"
$ws.Cells.Item(51,5).Value  = 50
$ws.Cells.Item(51,6).Value  = "E"
$ws.Cells.Item(51,7).Value  = "5. Software Testing KA"
$ws.Cells.Item(51,8).Value  = "5.2.3. Test Environment Set-up and Maintenance"
$ws.Cells.Item(51,9).Value  = "Para. 1"
$ws.Cells.Item(51,10).Value = "Line 9"

$k51 = $ws.Cells.Item(51,11)
$k51text = 'Since the terms "in vitro" and "in vivo" are never used in this document, it is unclear what they mean in the context of software engineering; does this mean that the testing environment can be either the actual environment or a simulated one (which seems redudant, since "simulated" is already given)?'
$k51.Value = $k51text
$run1 = $k51.Characters(206,6)
$run1.Font.Italic = $true
$run1.Font.Name = "Times New Roman"
$run1.Font.Size = 10
$run1.Font.ColorIndex = -4105
$run2 = $k51.Characters(230,9)
$run2.Font.Italic = $true
$run2.Font.Name = "Times New Roman"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105

$ws.Cells.Item(51,12).Value = 'Replace "in vitro or in vivo" with what that means in the context of software testing'

# ---- Row 52 --------------------------------------------------------
$ws.Cells.Item(52,1).Value  = "Samuel J. Crawford"
$ws.Cells.Item(52,2).Value  = "crawfs1@mcmaster.ca"
$ws.Cells.Item(52,3).Value  = "McMaster University"
$ws.Cells.Item(52,4).Value  = "Yes"
$ws.Cells.Item(52,5).Value  = 51
$ws.Cells.Item(52,6).Value  = "E"
$ws.Cells.Item(52,7).Value  = "5. Software Testing KA"
$ws.Cells.Item(52,8).Value  = "5.3. Staffing"
$ws.Cells.Item(52,9).Value  = "Paras. 1-2"
$ws.Cells.Item(52,10).Value = "Lines 7-8"
$ws.Cells.Item(52,11).Value = "A page break was added erroneously"
$ws.Cells.Item(52,12).Value = '"…to meet deadlines, and increase/reduce maintenance costs." (with no paragraph break)'

# ---------------------------------------------------------------------
# 2. Hyperlink the e-mail addresses in the two new rows (added in the
#    same order Excel assigns them: B52 first, then B51).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(52,2), "mailto:crawfs1@mcmaster.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(51,2), "mailto:crawfs1@mcmaster.ca") | Out-Null

# ---------------------------------------------------------------------
# 3. Row 62 changes from the old "last blank row before filler" style to
#    the regular in-range blank-row style (matches rows 53-61), since the
#    reviewed range grew by two rows.
# ---------------------------------------------------------------------
$ws.Range("G53").Copy($ws.Range("G62"))
$ws.Range("J53").Copy($ws.Range("J62"))
$ws.Rows.Item(62).RowHeight = 12.9

# ---------------------------------------------------------------------
# 4. A new filler row 130 appears at the bottom, matching row 129.
# ---------------------------------------------------------------------
$ws.Range("G129").Copy($ws.Range("G130"))
$ws.Range("J129").Copy($ws.Range("J130"))

# ---------------------------------------------------------------------
# 5. Data validation ranges grow along with the used range.
# ---------------------------------------------------------------------
$ws.Range("D2:D52").Validation.Delete()
$ws.Range("D2:D52").Validation.Add(3, 1, 1, "Yes,No")
$ws.Range("F2:F807").Validation.Delete()
$ws.Range("F2:F807").Validation.Add(3, 1, 1, "=Category")

# ---------------------------------------------------------------------
# 6. Selection ends up on the last edited cell, as it would after typing
#    the new rows in.
# ---------------------------------------------------------------------
$ws.Range("L52").Select()

Write-Output "applied SWEBOK testing-chapter comment rows"
